$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.256.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.322.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.896.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.312.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.297.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "442.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.23"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.521"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.470.27"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.194"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.73"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.851.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.790"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.43"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0667"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "328.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0273"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.15"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.20%  "
